$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 09:35"

# --- Update country stats: Singapur (row 30) ---
$ws.Range("B30").Value = 29364
$ws.Range("C30").Value = 570
$ws.Range("E30").Value = 18977

# --- Update country stats: Polonia (row 34) ---
$ws.Range("D34").Value = 8183
$ws.Range("E34").Value = 10137

# --- Update country stats: Chequia (row 52) ---
$ws.Range("B52").Value = 8683
$ws.Range("C52").Value = 36
$ws.Range("D52").Value = 5731
$ws.Range("E52").Value = 2649
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 303

# --- Re-order the small-countries block: move "Belice" up so it sits
#     right after "Fiyi" (ahead of "Nueva Caledonia" and "Santa Lucia"),
#     pushing those two down by one row, while keeping their own stats
#     attached to their own names. Rows 195-197 become:
#       195 Belice          (was row 197's data)
#       196 Nueva Caledonia (was row 195's data)
#       197 Santa Lucia     (was row 196's data)

$ws.Range("A195").Value = "Belice"
$ws.Range("B195").Value = 18
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 16
$ws.Range("E195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 2

$ws.Range("A196").Value = "Nueva Caledonia"
$ws.Range("B196").Value = 18
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 18
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("B197").Value = 18
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 18
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0
